# Apply metadata and element-table updates to the SSI Status StructureDefinition workbook.

$wb = $excel.ActiveWorkbook

# --- "Metadata" worksheet: update URL, Version, Date and Publisher values ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/ssi-status"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" worksheet: clear the stale Constraint(s) value on the root Extension row ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""

# The extension's fixed URL value mirrors the StructureDefinition URL above, keep them in sync
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/ssi-status"
